$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 60: keep the same feedback text, but store politeness_score (col B)
# as a real number (3) instead of text.
$ws.Range("B60").Value = 3

# Row 61 is new: it holds what used to be row 60's feedback text (id,
# source_file, text, issue_type) but with sentence_purpose changed to
# "DIS" and politeness_score stored as literal text "3" (not a number).
$ws.Range("A61").Value = "Ruilin"

# Build the text "3" without ever creating a numeric/quote-prefixed cell
# style: compute it as a formula result, then paste-as-values so the cell
# ends up a plain shared-string "3", matching t="inlineStr" in the source.
$ws.Range("ZZ1").Formula = '=TEXT(3,"0")'
$ws.Range("ZZ1").Copy()
$ws.Range("B61").PasteSpecial(-4163)
$ws.Range("ZZ1").Clear()

$ws.Range("C61").Value = "无"
$ws.Range("D61").Value = "DIS"
$ws.Range("E61").Value = "WRI"
$ws.Range("F61").Value = "1e0176d5-be35-49c3-adce-f7bfa3b6964a"
$ws.Range("G61").Value = "HksxTdiWz_annotated.xlsx"
$ws.Range("H61").Value = "In any case, this statement should be clarified."
